$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd($wb, $name) {
    $count = $wb.Worksheets.Count
    $last = $wb.Worksheets.Item($count)
    $newSheet = $wb.Worksheets.Add($null, $last)
    $newSheet.Name = $name
    return $newSheet
}

$neo4jUrlLabel = 'Neo4j_URL:'
$neo4jUrlValue = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$userNameLabel = 'User_name:'
$userNameValue = 'neo4j'
$pwdLabel = 'PWD:'
$pwdValue = 'icdcDBneo4j0'
$cypherLabel = 'Cypher:'
$cypherValue = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Lymphoma''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$outputLabel = 'Output:'
$outputValue = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC05_Canine_Filter_Diagnosis-Lymphoma_Neo4jData.xlsx'
$cypherEmptyMessage = 'Cypher query should not be an empty string'

$cypherOutputMessage = Add-SheetAtEnd $wb "CypherOutput_Message"
$cypherOutputMessage.Range("A1").Value = $neo4jUrlLabel
$cypherOutputMessage.Range("A2").Value = $neo4jUrlValue
$cypherOutputMessage.Range("A3").Value = $userNameLabel
$cypherOutputMessage.Range("A4").Value = $userNameValue
$cypherOutputMessage.Range("A5").Value = $pwdLabel
$cypherOutputMessage.Range("A6").Value = $pwdValue
$cypherOutputMessage.Range("A7").Value = $cypherLabel
$cypherOutputMessage.Range("A8").Value = $cypherValue
$cypherOutputMessage.Range("A9").Value = $outputLabel
$cypherOutputMessage.Range("A10").Value = $outputValue

$statOutput = Add-SheetAtEnd $wb "StatOutput"

$statOutputMessage = Add-SheetAtEnd $wb "StatOutput_Message"
$statOutputMessage.Range("A1").Value = $neo4jUrlLabel
$statOutputMessage.Range("A2").Value = $neo4jUrlValue
$statOutputMessage.Range("A3").Value = $userNameLabel
$statOutputMessage.Range("A4").Value = $userNameValue
$statOutputMessage.Range("A5").Value = $pwdLabel
$statOutputMessage.Range("A6").Value = $pwdValue
$statOutputMessage.Range("A7").Value = $cypherLabel
$statOutputMessage.Range("A8").Value = $cypherValue
$statOutputMessage.Range("A9").Value = $outputLabel
$statOutputMessage.Range("A10").Value = $outputValue
$statOutputMessage.Range("A11").Value = $cypherEmptyMessage
$statOutputMessage.Range("A12").Value = $neo4jUrlLabel
$statOutputMessage.Range("A13").Value = $neo4jUrlValue
$statOutputMessage.Range("A14").Value = $userNameLabel
$statOutputMessage.Range("A15").Value = $userNameValue
$statOutputMessage.Range("A16").Value = $pwdLabel
$statOutputMessage.Range("A17").Value = $pwdValue
$statOutputMessage.Range("A18").Value = $cypherLabel
$statOutputMessage.Range("A19").Value = "'"
$statOutputMessage.Range("A19").ClearFormats()
$statOutputMessage.Range("A20").Value = $outputLabel
$statOutputMessage.Range("A21").Value = $outputValue

$wb.Worksheets.Item("CypherOutput").Activate()

